{"js": "// Apply the dated-worksheet update: bump the header date by one day and\n// replace each division-problem cell with its new problem, in document order.\n// Some search strings (e.g. \"73\u00f73=\") occur more than once, each time mapping\n// to a different replacement, so we resolve every occurrence by its index in\n// the document rather than doing a blind find-and-replace-all.\n\nconst body = context.document.body;\n\n// [searchText, [replacement-for-1st-hit, replacement-for-2nd-hit, ...]]\nconst replacements = [\n  [\"2025-08-24 Sunday\", [\"2025-08-25 Monday\"]],\n  [\"65\u00f74=\", [\"37\u00f73=\"]],\n  [\"46\u00f79=\", [\"53\u00f74=\"]],\n  [\"33\u00f75=\", [\"29\u00f72=\"]],\n  [\"73\u00f73=\", [\"43\u00f77=\", \"95\u00f78=\", \"50\u00f75=\"]],\n  [\"80\u00f78=\", [\"98\u00f77=\"]],\n  [\"79\u00f78=\", [\"84\u00f72=\"]],\n  [\"35\u00f79=\", [\"65\u00f72=\"]],\n  [\"32\u00f75=\", [\"26\u00f79=\"]],\n  [\"31\u00f75=\", [\"41\u00f79=\"]],\n  [\"21\u00f75=\", [\"38\u00f72=\"]],\n  [\"99\u00f74=\", [\"12\u00f72=\"]],\n  [\"17\u00f74=\", [\"63\u00f76=\"]],\n  [\"59\u00f78=\", [\"62\u00f76=\"]],\n  [\"53\u00f79=\", [\"35\u00f74=\"]],\n  [\"56\u00f79=\", [\"65\u00f73=\"]],\n  [\"90\u00f78=\", [\"56\u00f72=\"]],\n  [\"36\u00f72=\", [\"17\u00f78=\"]],\n  [\"16\u00f74=\", [\"95\u00f72=\"]],\n  [\"30\u00f78=\", [\"66\u00f78=\"]],\n  [\"33\u00f76=\", [\"18\u00f78=\"]],\n  [\"78\u00f74=\", [\"72\u00f75=\"]],\n  [\"69\u00f72=\", [\"70\u00f72=\"]],\n  [\"76\u00f74=\", [\"58\u00f77=\"]],\n];\n\n// Search for every distinct old value once, load the matching ranges, then\n// apply the ordered replacements against the hits (which come back in\n// document order).\nconst searchResults = [];\nfor (const [searchText] of replacements) {\n  const found = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  searchResults.push(found);\n}\n\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newValues] = replacements[i];\n  const items = searchResults[i].items;\n  for (let j = 0; j < newValues.length; j++) {\n    items[j].insertText(newValues[j], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply the dated-worksheet update: bump the header date by one day and\n# replace each division-problem cell with its new problem.\n# Some search strings (e.g. \"73\u00f73=\") occur more than once in the document,\n# each occurrence mapping to a different replacement value, so each old\n# value is searched for sequentially (moving the search range forward after\n# every hit) and the Nth hit is replaced with the Nth new value in the list.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2025-08-24 Sunday\"; New = @(\"2025-08-25 Monday\") },\n    @{ Old = \"65\u00f74=\"; New = @(\"37\u00f73=\") },\n    @{ Old = \"46\u00f79=\"; New = @(\"53\u00f74=\") },\n    @{ Old = \"33\u00f75=\"; New = @(\"29\u00f72=\") },\n    @{ Old = \"73\u00f73=\"; New = @(\"43\u00f77=\", \"95\u00f78=\", \"50\u00f75=\") },\n    @{ Old = \"80\u00f78=\"; New = @(\"98\u00f77=\") },\n    @{ Old = \"79\u00f78=\"; New = @(\"84\u00f72=\") },\n    @{ Old = \"35\u00f79=\"; New = @(\"65\u00f72=\") },\n    @{ Old = \"32\u00f75=\"; New = @(\"26\u00f79=\") },\n    @{ Old = \"31\u00f75=\"; New = @(\"41\u00f79=\") },\n    @{ Old = \"21\u00f75=\"; New = @(\"38\u00f72=\") },\n    @{ Old = \"99\u00f74=\"; New = @(\"12\u00f72=\") },\n    @{ Old = \"17\u00f74=\"; New = @(\"63\u00f76=\") },\n    @{ Old = \"59\u00f78=\"; New = @(\"62\u00f76=\") },\n    @{ Old = \"53\u00f79=\"; New = @(\"35\u00f74=\") },\n    @{ Old = \"56\u00f79=\"; New = @(\"65\u00f73=\") },\n    @{ Old = \"90\u00f78=\"; New = @(\"56\u00f72=\") },\n    @{ Old = \"36\u00f72=\"; New = @(\"17\u00f78=\") },\n    @{ Old = \"16\u00f74=\"; New = @(\"95\u00f72=\") },\n    @{ Old = \"30\u00f78=\"; New = @(\"66\u00f78=\") },\n    @{ Old = \"33\u00f76=\"; New = @(\"18\u00f78=\") },\n    @{ Old = \"78\u00f74=\"; New = @(\"72\u00f75=\") },\n    @{ Old = \"69\u00f72=\"; New = @(\"70\u00f72=\") },\n    @{ Old = \"76\u00f74=\"; New = @(\"58\u00f77=\") }\n)\n\nforeach ($item in $replacements) {\n    $rng = $d.Content\n    $rng.Start = 0\n    $i = 0\n    while ($rng.Find.Execute($item.Old)) {\n        $rng.Text = $item.New[$i]\n        $i = $i + 1\n        $rng.Collapse(0)\n    }\n}\n"}
